$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand the table (ListObject) from A1:B29 to A1:B33
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B33"))

# Column A: header "ID" then sequential numbers 1..32
$ws.Cells.Item(1, 1).Value = "ID"
for ($i = 2; $i -le 33; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
}

# Column B: header "İsim" then course names
$ws.Cells.Item(1, 2).Value = "İsim"
$ws.Cells.Item(2, 2).Value = "Bilgisayar Mühendisliğine Giriş"
$ws.Cells.Item(3, 2).Value = "Bilgisayar Programlama 1"
$ws.Cells.Item(4, 2).Value = "Fizik 1"
$ws.Cells.Item(5, 2).Value = "Matematik 1"
$ws.Cells.Item(6, 2).Value = "Diferansiyel Denklemler"
$ws.Cells.Item(7, 2).Value = "Elektrik Mühendisliğinin Temelleri"
$ws.Cells.Item(8, 2).Value = "Nesne Yönelimli Programlama"
$ws.Cells.Item(9, 2).Value = "Olasılık ve İstatistik"
$ws.Cells.Item(10, 2).Value = "Bilişim Etiği ve Hukuku"
$ws.Cells.Item(11, 2).Value = "İşletim Sistemleri"
$ws.Cells.Item(12, 2).Value = "Mantıksal Devre Tasarımı"
$ws.Cells.Item(13, 2).Value = "Sayısal Analiz"
$ws.Cells.Item(14, 2).Value = "Siber Güvenlik"
$ws.Cells.Item(15, 2).Value = "Veritabanı Yönetimi"
$ws.Cells.Item(16, 2).Value = "Makine Öğrenmesi"
$ws.Cells.Item(17, 2).Value = "Mobil Programlama"
$ws.Cells.Item(18, 2).Value = "Sayısal Görüntü İşleme"
$ws.Cells.Item(19, 2).Value = "Yapay Zekaya Giriş"
$ws.Cells.Item(20, 2).Value = "Endüstri Mühendisliğine Giriş"
$ws.Cells.Item(21, 2).Value = "Lineer Cebir"
$ws.Cells.Item(22, 2).Value = "Genel Muhasebe"
$ws.Cells.Item(23, 2).Value = "Malzeme Bilimi"
$ws.Cells.Item(24, 2).Value = "Mühendislik Ekonomisi"
$ws.Cells.Item(25, 2).Value = "Stok ve Envanter Yönetimi"
$ws.Cells.Item(26, 2).Value = "Yöneylem Araştırması 1"
$ws.Cells.Item(27, 2).Value = "Ergonomi"
$ws.Cells.Item(28, 2).Value = "Benzetim"
$ws.Cells.Item(29, 2).Value = "Sistem Analizi ve Tasarımı"
$ws.Cells.Item(30, 2).Value = "Girişimcilik"
$ws.Cells.Item(31, 2).Value = "Üretim Planlama ve Kontrol"
$ws.Cells.Item(32, 2).Value = "Çizelgeleme"
$ws.Cells.Item(33, 2).Value = "Çok Ölçütlü Karar Verme"

# Update selection to match target state
$ws.Range("G24").Select()
